# Auto-generated edit script: updates crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "23.066.54"
Set-TextValue "E2" "  -0.39%  "
Set-TextValue "D3" "1.592.81"
Set-TextValue "E3" "  -0.41%  "
Set-TextValue "D4" "1.001"
Set-TextValue "E4" "  -0.12%  "
Set-TextValue "E5" "  +0.00%  "
Set-TextValue "D6" "302.01"
Set-TextValue "E6" "  +0.24%  "
Set-TextValue "D7" "0.3775"
Set-TextValue "E7" "  +0.16%  "
Set-TextValue "D8" "0.3608"
Set-TextValue "E8" "  -1.31%  "
Set-TextValue "D9" "50.61"
Set-TextValue "E9" "  +6.58%  "
Set-TextValue "D10" "1.001"
Set-TextValue "E10" "  -0.10%  "
Set-TextValue "D11" "1.233"
Set-TextValue "E11" "  -3.46%  "
Set-TextValue "D12" "0.08105"
Set-TextValue "E12" "  +0.33%  "
Set-TextValue "D13" "22.16"
Set-TextValue "E13" "  -3.36%  "
Set-TextValue "D14" "6.514"
Set-TextValue "E14" "  -1.85%  "
Set-TextValue "D15" "7.286"
Set-TextValue "E15" "  -3.73%  "
Set-TextValue "D16" "0.00001232"
Set-TextValue "E16" "  -2.64%  "
Set-TextValue "D17" "1.590.31"
Set-TextValue "E17" "  -0.39%  "
Set-TextValue "D18" "92.50"
Set-TextValue "E18" "  +0.93%  "
Set-TextValue "D19" "0.06833"
Set-TextValue "E19" "  +0.40%  "
Set-TextValue "D20" "18.06"
Set-TextValue "E20" "  -2.05%  "
Set-TextValue "D21" "6.485"
Set-TextValue "E21" "  -1.74%  "
Set-TextValue "D22" "1.003"
Set-TextValue "E22" "  -0.02%  "
Set-TextValue "E23" "  -0.93%  "
Set-TextValue "D24" "23.069.88"
Set-TextValue "E24" "  -0.36%  "
Set-TextValue "E25" "  +0.65%  "
Set-TextValue "D26" "2.813"
Set-TextValue "E26" "  -3.33%  "
Set-TextValue "D27" "21.01"
Set-TextValue "E27" "  -0.47%  "
Set-TextValue "D28" "148.73"
Set-TextValue "E28" "  -1.52%  "
Set-TextValue "D29" "5.236"
Set-TextValue "E29" "  -0.19%  "
Set-TextValue "D30" "133.61"
Set-TextValue "E30" "  +1.08%  "
Set-TextValue "D31" "2.375"
Set-TextValue "E31" "  -2.60%  "
Set-TextValue "D32" "6.669"
Set-TextValue "E32" "  -6.39%  "
Set-TextValue "D33" "1.764.65"
Set-TextValue "E33" "  -0.43%  "
Set-TextValue "D34" "0.9526"
Set-TextValue "E34" "  -2.40%  "
Set-TextValue "D35" "0.07440"
Set-TextValue "E35" "  -3.88%  "
Set-TextValue "D36" "10.11"
Set-TextValue "E36" "  +0.06%  "
Set-TextValue "D37" "0.02691"
Set-TextValue "E37" "  -3.31%  "
Set-TextValue "B38" "InternetComputer(DFINITY)"
Set-TextValue "C38" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D38" "6.114"
Set-TextValue "E38" "  -3.20%  "
Set-TextValue "B39" "Stellar"
Set-TextValue "C39" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D39" "0.08810"
Set-TextValue "E39" "  -0.67%  "
Set-TextValue "D40" "0.2496"
Set-TextValue "E40" "  -2.01%  "
Set-TextValue "D41" "1.359"
Set-TextValue "E41" "  -3.28%  "
Set-TextValue "D42" "0.6975"
Set-TextValue "E42" "  -2.64%  "
Set-TextValue "D43" "12.26"
Set-TextValue "E43" "  -3.98%  "
Set-TextValue "E44" "  -6.21%  "
Set-TextValue "D45" "0.6510"
Set-TextValue "E45" "  -1.98%  "
Set-TextValue "D46" "4.010"
Set-TextValue "E46" "  +1.12%  "
Set-TextValue "D47" "2.269"
Set-TextValue "E47" "  -1.83%  "
Set-TextValue "D48" "131.65"
Set-TextValue "E48" "  +0.08%  "
Set-TextValue "D49" "0.07908"
Set-TextValue "E49" "  -1.02%  "
Set-TextValue "D50" "1.209"
Set-TextValue "E50" "  +3.27%  "
Set-TextValue "D51" "1.220"
Set-TextValue "E51" "  +4.44%  "
